# Update cryptocurrency "Price" (column D) values on the active sheet
# as re-scraped by the GitHub Actions job on Wed Dec 21 13:57:38 UTC 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "249.18"
    3  = "22.68"
    4  = "5.264"
    5  = "0.05694"
    6  = "3.408"
    7  = "6.336"
    8  = "0.8052"
    9  = "0.8999"
    10 = "0.1417"
    11 = "0.07437"
    12 = "0.03100"
    14 = "0.09386"
    15 = "3.873"
    16 = "0.001593"
    17 = "0.04755"
    19 = "0.0005812"
    20 = "0.006456"
    21 = "0.004986"
    22 = "0.0009997"
    24 = "3.694"
    25 = "2.196"
    27 = "0.1355"
    40 = "0.03975"
    41 = "0.006688"
    42 = "0.1071"
    43 = "0.002729"
    44 = "0.007726"
    47 = "0.4993"
    48 = "0.2051"
    49 = "0.00002101"
    50 = "0.01011"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    # Force the new value to be stored as text (preserving exact digits,
    # e.g. trailing zeros) rather than being auto-converted to a number.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    # Restore the cell's original (unstyled) appearance so only the
    # textual content changes, matching the source data's plain format.
    $cell.Style = "Normal"
}
